$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Combined Ratio", 91, "Definity", "Q1 2024"),
    @("Combined Ratio", 94.09999999999999, "Definity", "Q1 2025"),
    @("Claims Ratio", 59.9, "Definity", "Q1 2025"),
    @("Claims Ratio", 55.2, "Definity", "Q1 2024"),
    @("Core Claim Ratio", 51.3, "Definity", "Q1 2024"),
    @("Core Claim Ratio", 52, "Definity", "Q1 2025"),
    @("CAT Loss Ratio", 11.8, "Definity", "Q1 2025"),
    @("CAT Loss Ratio", 5.9, "Definity", "Q1 2024"),
    @("Expense Ratio", 35.8, "Definity", "Q1 2024"),
    @("Expense Ratio", 34.2, "Definity", "Q1 2025"),
    @("PYD Ratio", -3.9, "Definity", "Q1 2025"),
    @("PYD Ratio", -2, "Definity", "Q1 2024"),
    @("Gross Written Premium", 236.5, "Definity", "Q1 2024"),
    @("Gross Written Premium", 255, "Definity", "Q1 2025"),
    @("Underwriting Income", 16.8, "Definity", "Q1 2025"),
    @("Underwriting Income", 23.5, "Definity", "Q1 2024"),
    @("ROE", 9.5, "Definity", "Q1 2024"),
    @("ROE", 10.3, "Definity", "Q1 2025"),
    @("Combined Ratio", 88.90000000000001, "Intact", "Q1 2025"),
    @("Combined Ratio", 82.5, "Intact", "Q1 2024"),
    @("Claims Ratio", 46.4, "Intact", "Q1 2024"),
    @("Claims Ratio", 55.6, "Intact", "Q1 2025"),
    @("Core Claim Ratio", 53.7, "Intact", "Q1 2025"),
    @("Core Claim Ratio", 51, "Intact", "Q1 2024"),
    @("CAT Loss Ratio", 0, "Intact", "Q1 2024"),
    @("CAT Loss Ratio", 7.5, "Intact", "Q1 2025"),
    @("Expense Ratio", 33.3, "Intact", "Q1 2025"),
    @("Expense Ratio", 36.1, "Intact", "Q1 2024"),
    @("PYD Ratio", -4.6, "Intact", "Q1 2024"),
    @("PYD Ratio", -5.6, "Intact", "Q1 2025"),
    @("Gross Written Premium", 903, "Intact", "Q1 2025"),
    @("Gross Written Premium", 828, "Intact", "Q1 2024"),
    @("Underwriting Income", 166, "Intact", "Q1 2024"),
    @("Underwriting Income", 113, "Intact", "Q1 2025"),
    @("ROE", 13.7, "Intact", "Q1 2025"),
    @("ROE", 10.6, "Intact", "Q1 2024")
)

$startRow = 202
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
